# Apply text/wording tweaks to the Language workbook ("en" sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# --- Text corrections (Column B holds the localized value for the Key in Column A) ---
# (Assignment order matches the order new shared-string entries were authored in,
#  so the rebuilt sharedStrings table lines up with the target.)

# intro_0_1 (row 81): "fine capable frogs analyze" -> "fine and capable frogs to analyze"
$ws.Range("B81").Value = "In this game, you will be helping these fine and capable frogs to analyze a potentially habitable planet by identifying the shapes scattered across the land."

# ellipse_desc (row 69): "two axis" -> "two axes"
# NOTE: this workbook stores "\n" as a literal two-character escape sequence
# (backslash + n) inside the string value, not an actual line break -- so we
# must NOT use PowerShell's `n (real newline) escape here.
$ws.Range("B69").Value = '· Has a center point.\n\n· Roundness defined by two axes.'

# poly_desc (row 71): comma -> period, new sentence
$ws.Range("B71").Value = "· Formed by three or more straight lines (sides) connected in a loop. Points are plotted on a plane."

# level_intro_4_0 (row 93): "the ones based on" -> "some more based on"
$ws.Range("B93").Value = "Now that we've seen some triangles based on their angle values, let's take a look at some more based on their side values."

# level_intro_5_0 (row 94): "angles, and side lengths" -> "angles AND side lengths"
$ws.Range("B94").Value = "For this level, we will be categorizing triangles based on their angles AND side lengths. These triangles will have more than one category that fit their attributes."

# --- Update the view state to match the saved selection in the target file ---
$excel.Goto($ws.Range("B94"), $true)
